$wb = $excel.ActiveWorkbook

# --- Column order shared by the "Confirmados" (sheet1) and "Mortes" (sheet2)
# tables: B..AB = Acre, Alagoas, Amapa, Amazonas, Bahia, Ceara,
# Distrito Federal, Espirito Santo, Goias, Maranhao, Mato Grosso,
# Mato Grosso do Sul, Minas Gerais, Para, Paraiba, Parana, Pernambuco,
# Piaui, Rio de Janeiro, Rio Grande do Norte, Rio Grande do Sul, Rondonia,
# Roraima, Santa Catarina, Sao Paulo, Sergipe, Tocantins

$confirmados = @(90,50,242,1275,723,1800,638,430,233,445,134,113,815,270,111,756,1154,50,3231,339,664,42,83,777,8895,44,26)
$mortes      = @(3,3,5,71,22,91,15,14,15,27,4,4,23,15,13,31,102,8,188,17,16,2,3,24,608,4,0)

$newRow = 50

# Build the new date label ("2020-04-13") as plain text in a scratch cell via
# a text formula, then paste-special the computed value into the new date
# cell. This avoids Excel's automatic "looks like a date" conversion that a
# direct Range.Value assignment of an ISO date string would trigger, while
# leaving no stray number-format/style behind.
$scratch = $wb.Worksheets.Item(1).Range("AZ1000")
$scratch.Formula = "=""2020-04-13"""
$scratch.Copy()
$wb.Worksheets.Item(1).Range("A" + $newRow).PasteSpecial(-4163)
$scratch.ClearContents()

$scratch2 = $wb.Worksheets.Item(2).Range("AZ1000")
$scratch2.Formula = "=""2020-04-13"""
$scratch2.Copy()
$wb.Worksheets.Item(2).Range("A" + $newRow).PasteSpecial(-4163)
$scratch2.ClearContents()

# Fill in the numeric columns for both sheets.
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
for ($i = 0; $i -lt $confirmados.Length; $i++) {
    $ws1.Cells.Item($newRow, 2 + $i).Value = $confirmados[$i]
    $ws2.Cells.Item($newRow, 2 + $i).Value = $mortes[$i]
}
